$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 101 (shifts existing rows 101-109 down to 102-110),
# matching the author-list insertion of "Teruel-Pardo, S." (IFIC / Paterna)
# right before "Toledo, J.F.".
$ws.Rows.Item(101).Insert()

$ws.Range("A101").Value = 'Teruel-Pardo'
$ws.Range("B101").Value = 'S.'
$ws.Range("E101").Value = 'Instituto de F\''isica Corpuscular (IFIC), CSIC \& Universitat de Val\`encia, Calle Catedr\''atico Jos\''e Beltr\''an, 2 '
$ws.Range("F101").Value = ' Paterna, E-46980, Spain'

$excel.ActiveWindow.ScrollRow = 69
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D101").Select()
